# Generate Report for Handoff
#
# Adds two new source files (d004172b-...md and d4bd82b5-...md) to the
# localization-status report. Each gets a row in the "Overview" sheet and
# in the per-locale "zh-cn" / "de-de" sheets, inserted just above the
# existing ".localization-config" bookkeeping row (which shifts down).

$wb = $excel.ActiveWorkbook

$newFile1 = "d004172b-70b9-449d-89e0-0618d317ef77.md"
$newFile2 = "d4bd82b5-bd23-4684-9f56-4a2aaaa8fde9.md"

$xlf1_zh = "d004172b-70b9-449d-89e0-0618d317ef77.5a9611d558c06e4b669ac631e927ffafb9b92ffe.zh-cn.xlf"
$xlf2_zh = "d4bd82b5-bd23-4684-9f56-4a2aaaa8fde9.2c98e03ce57e7141d52c7e6dcde628f1b36de271.zh-cn.xlf"
$xlf1_de = "d004172b-70b9-449d-89e0-0618d317ef77.5a9611d558c06e4b669ac631e927ffafb9b92ffe.de-de.xlf"
$xlf2_de = "d4bd82b5-bd23-4684-9f56-4a2aaaa8fde9.2c98e03ce57e7141d52c7e6dcde628f1b36de271.de-de.xlf"

$handoffDt_zh = "2016-02-29 11:52:54"
$handoffDt_de = "2016-02-29 11:53:04"

$commitBase = "ac323ca43e30c67d51654c80e6c0edf635c7e1e4"
$zhCommit = "e027c3a8457d849a9c8493c560c963bf6e3ecbc7"
$deCommit = "2e63910c667b0294533447ab7cd46f71aed254f9"

# ---------------------------------------------------------------------
# Sheet "Overview" (File Name / zh-cn / de-de summary)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Make room: push the ".localization-config" row (currently row 4) down
# to row 6, inheriting its style onto two freshly inserted blank rows.
$ws1.Rows.Item(4).Insert()
$ws1.Rows.Item(4).Insert()

$ws1.Cells.Item(4,1).Value = $newFile1
$ws1.Cells.Item(4,2).Value = "Ready for handoff"
$ws1.Cells.Item(4,3).Value = "Ready for handoff"

$ws1.Cells.Item(5,1).Value = $newFile2
$ws1.Cells.Item(5,2).Value = "Ready for handoff"
$ws1.Cells.Item(5,3).Value = "Ready for handoff"

# Row 6 already carries the correct ".localization-config" / "Not to be
# localized" content (and style) courtesy of the row insert above.

# Rebuild hyperlinks in final top-to-bottom order so relationship ids
# come out sequential and every display text matches its row.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/a2853b7d-4952-4a93-84bd-6e89444f2370.md", "", "", "a2853b7d-4952-4a93-84bd-6e89444f2370.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/fc3b4349-984b-4737-8f7b-eb5f6be56167.md", "", "", "fc3b4349-984b-4737-8f7b-eb5f6be56167.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile1", "", "", $newFile1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile2", "", "", $newFile2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(4).Insert()
$ws2.Rows.Item(4).Insert()

$ws2.Cells.Item(4,1).Value = $newFile1
$ws2.Cells.Item(4,2).Value = "Ready for handoff"
$ws2.Cells.Item(4,3).Value = $xlf1_zh
$ws2.Cells.Item(4,4).Value = $handoffDt_zh
$ws2.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,8).Value = "Include"

$ws2.Cells.Item(5,1).Value = $newFile2
$ws2.Cells.Item(5,2).Value = "Ready for handoff"
$ws2.Cells.Item(5,3).Value = $xlf2_zh
$ws2.Cells.Item(5,4).Value = $handoffDt_zh
$ws2.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5,8).Value = "Include"

# Row 6 already carries the correct ".localization-config" data/style.

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/a2853b7d-4952-4a93-84bd-6e89444f2370.md", "", "", "a2853b7d-4952-4a93-84bd-6e89444f2370.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a2853b7d-4952-4a93-84bd-6e89444f2370.76bda225dedd27874e1c95f85cee2502b878d3e8.zh-cn.xlf", "", "", "a2853b7d-4952-4a93-84bd-6e89444f2370.76bda225dedd27874e1c95f85cee2502b878d3e8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/fc3b4349-984b-4737-8f7b-eb5f6be56167.md", "", "", "fc3b4349-984b-4737-8f7b-eb5f6be56167.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc3b4349-984b-4737-8f7b-eb5f6be56167.71ffe595d5cdfcf85b18ed222731b2bec0eba958.zh-cn.xlf", "", "", "fc3b4349-984b-4737-8f7b-eb5f6be56167.71ffe595d5cdfcf85b18ed222731b2bec0eba958.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile1", "", "", $newFile1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlf1_zh", "", "", $xlf1_zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile2", "", "", $newFile2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlf2_zh", "", "", $xlf2_zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(4).Insert()
$ws3.Rows.Item(4).Insert()

$ws3.Cells.Item(4,1).Value = $newFile1
$ws3.Cells.Item(4,2).Value = "Ready for handoff"
$ws3.Cells.Item(4,3).Value = $xlf1_de
$ws3.Cells.Item(4,4).Value = $handoffDt_de
$ws3.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,8).Value = "Include"

$ws3.Cells.Item(5,1).Value = $newFile2
$ws3.Cells.Item(5,2).Value = "Ready for handoff"
$ws3.Cells.Item(5,3).Value = $xlf2_de
$ws3.Cells.Item(5,4).Value = $handoffDt_de
$ws3.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5,8).Value = "Include"

# Row 6 already carries the correct ".localization-config" data/style.

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/a2853b7d-4952-4a93-84bd-6e89444f2370.md", "", "", "a2853b7d-4952-4a93-84bd-6e89444f2370.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a2853b7d-4952-4a93-84bd-6e89444f2370.76bda225dedd27874e1c95f85cee2502b878d3e8.de-de.xlf", "", "", "a2853b7d-4952-4a93-84bd-6e89444f2370.76bda225dedd27874e1c95f85cee2502b878d3e8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/fc3b4349-984b-4737-8f7b-eb5f6be56167.md", "", "", "fc3b4349-984b-4737-8f7b-eb5f6be56167.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc3b4349-984b-4737-8f7b-eb5f6be56167.71ffe595d5cdfcf85b18ed222731b2bec0eba958.de-de.xlf", "", "", "fc3b4349-984b-4737-8f7b-eb5f6be56167.71ffe595d5cdfcf85b18ed222731b2bec0eba958.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile1", "", "", $newFile1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlf1_de", "", "", $xlf1_de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/e2e/$newFile2", "", "", $newFile2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlf2_de", "", "", $xlf2_de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commitBase/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Report updated: added $newFile1 and $newFile2 rows to Overview, zh-cn, de-de sheets"
